$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.036.23"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -1.52%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.859.28"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -2.82%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.26%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "232.68"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -2.75%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.31%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4662"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -2.49%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2805"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -1.51%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06528"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -2.75%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.00"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +5.42%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07740"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +0.43%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "95.85"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -6.70%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.869.09"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -2.57%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.039"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -3.49%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6654"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "280.67"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +3.12%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.077.66"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -1.46%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.Style = "Normal"
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.121.01"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -1.75%  "
$c.Style = "Normal"
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.Style = "Normal"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.47"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.317"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -2.10%  "
$c.Style = "Normal"
$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.Style = "Normal"
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.21%  "
$c.Style = "Normal"
$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.Style = "Normal"
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.000007196"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -3.76%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.104"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -3.42%  "
$c.Style = "Normal"
$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.Style = "Normal"
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "166.61"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.Style = "Normal"
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.259"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -1.61%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.87"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -2.61%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.937"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -6.35%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.363"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -1.43%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09625"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -4.24%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.343"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -5.77%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.460"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.49%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.063"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -4.19%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04635"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -2.13%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6947"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -4.82%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.079"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -3.01%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.705"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.44%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01845"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -4.20%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.356"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +0.94%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.502"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -4.23%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "70.93"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -5.54%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8562"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -0.26%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.924"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -2.44%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +0.34%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "102.71"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -2.16%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4142"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -3.00%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "975.70"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +5.93%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.139"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -4.00%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "9.054"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +2.43%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "33.64"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -3.55%  "
$c.Style = "Normal"
